$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B and C are treated as plain text so values like
# "$61.00" or "01/20/2026" are not auto-converted into numbers/dates.
$ws.Range("B2:C21").NumberFormat = "@"

$data = @(
    @{ Row = 2;  B = "$61.00";  C = "01/20/2026" },
    @{ Row = 3;  B = "$31.50";  C = "01/19/2026" },
    @{ Row = 4;  B = "$185.00"; C = "01/19/2026" },
    @{ Row = 5;  B = "$129.00"; C = "01/19/2026" },
    @{ Row = 6;  B = "$110.00"; C = "01/19/2026" },
    @{ Row = 7;  B = "$48.00";  C = "01/13/2026" },
    @{ Row = 8;  B = "$21.00";  C = "01/13/2026" },
    @{ Row = 9;  B = "$171.00"; C = "01/13/2026" },
    @{ Row = 10; B = "$21.00";  C = "01/12/2026" },
    @{ Row = 11; B = "$21.00";  C = "01/12/2026" },
    @{ Row = 12; B = "$21.00";  C = "01/12/2026" },
    @{ Row = 13; B = "$69.50";  C = "01/12/2026" },
    @{ Row = 14; B = "$69.50";  C = "01/12/2026" },
    @{ Row = 15; B = "$31.50";  C = "01/12/2026" },
    @{ Row = 16; B = "$195.50"; C = "01/12/2026" },
    @{ Row = 17; B = "$106.50"; C = "01/12/2026" },
    @{ Row = 18; B = "$77.00";  C = "01/12/2026" },
    @{ Row = 19; B = "$21.00";  C = "01/12/2026" },
    @{ Row = 20; B = "$21.00";  C = "01/12/2026" },
    @{ Row = 21; B = "$21.00";  C = "01/12/2026" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
}
